$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New id / speaker_variant values for rows 2-19 (column B = id, column C = speaker_variant)
$data = @{
    2  = @("#tjemke-uyt:", "Tjemke uyt:")
    3  = @("#gerb", "Gerb")
    4  = @("#ias", "Ias")
    5  = @("#schout", "Schout")
    6  = @("#droncke-ger", "Droncke Ger")
    7  = @("#rem", "Rem")
    8  = @("#col", "Col")
    9  = @("#ian", "Ian")
    10 = @("#ron", "Ron")
    11 = @("#jan", "Jan")
    12 = @("#hy-haelt-refereyntjes-van-sijn-gordelen-leest", "Hy haelt refereyntjes van sijn gordelen leest")
    13 = @("#rons", "Rons")
    14 = @("#ger", "Ger")
    15 = @("#knip", "Knip")
    16 = @("#adv", "Adv")
    17 = @("#flor", "Flor")
    18 = @("#tjem", "Tjem")
    19 = @("#remmert-lubbertsz", "Remmert Lubbertsz")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    # Clear the "is_prefered" column (D) for every data row
    $ws.Cells.Item($row, 4).Value = ""
}
